# Update 15th March 2025 2328 Hours
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of Daily Code Count tracking data (rows 39-46)
$ws.Range("C39").Value = "Permutations II"
$ws.Range("E39").Value = "LeetCode"

$ws.Range("C40").Value = "Permutations II"
$ws.Range("E40").Value = "Bosscoder Academy"

$ws.Range("C41").Value = "Subsets II"
$ws.Range("E41").Value = "Leetcode"

$ws.Range("C42").Value = "Combination Sum (Solution 1)"
$ws.Range("E42").Value = "Leetcode"

$ws.Range("C43").Value = "Combination Sum II (Solution 1)"
$ws.Range("E43").Value = "LeetCode"

$ws.Range("B44").Value = 14
$ws.Range("C44").Value = "Count of Smaller Numbers After Self"
$ws.Range("E44").Value = "Bosscoder Academy"

$ws.Range("B45").Value = 15
$ws.Range("C45").Value = "Combination Sum (Solution 2)"
$ws.Range("E45").Value = "Leetcode"

$ws.Range("C46").Value = "Merge Sorted Array"
$ws.Range("E46").Value = "LeetCode"

# Scroll / selection state like the saved workbook
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("C46").Select()
